$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) must remain Text (the source feed always stores
# prices as text, even when they look numeric, e.g. "585.16"). Force the
# cell to Text format before assigning so Excel does not coerce it to a
# number, then clear the format again so the cell keeps the same (default,
# unstyled) appearance it had before the edit.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.907.55'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.94%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.503.90'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.71%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.16'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.56'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.65%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.504.65'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.68%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.484'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.31%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.123'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.09%  '
$ws.Range("E11").Value = '  -0.73%  '
$ws.Range("E12").Value = '  -2.83%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.105.35'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.41'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.88%  '
$ws.Range("E15").Value = '  +1.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.513.73'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.35%  '
$ws.Range("E17").Value = '  -2.31%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.960.23'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.75%  '
$ws.Range("E19").Value = '  -3.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.90'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.59'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.69%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '382.29'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.570'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.72%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.651.12'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.75'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.32%  '
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.65'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.78%  '
$ws.Range("E28").Value = '  +1.77%  '
$ws.Range("E29").Value = '  -2.35%  '
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.23%  '
$ws.Range("B31").Value = 'RenderToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.44'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -3.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.36'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.16%  '
$ws.Range("E33").Value = '  -2.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.518.84'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.47%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.45'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.84%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.144'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.32'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.23%  '
$ws.Range("E39").Value = '  -0.84%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.87'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '160.59'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -4.52%  '
$ws.Range("E42").Value = '  -3.02%  '
$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.811'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.55%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '26.43'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.92%  '
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '41.57'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.30%  '
$ws.Range("E47").Value = '  -4.61%  '
$ws.Range("E48").Value = '  -1.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.60'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -3.94%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.479.43'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.77'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.14%  '
